$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in rows 95 and 96 (column B) ---
$ws.Range("B95").Value = 129.949019653322
$ws.Range("B96").Value = 128.831398913345

# --- Append a new row 97 with a new monthly data point ---
# Column A: date (first of December 2024). Copy the date formatting from the
# cell above first (so it picks up the same date-number-format style), then
# overwrite the value.
$ws.Range("A96").Copy($ws.Range("A97"))
$ws.Range("A97").Value = 45627

# Columns B & C: plain numeric values
$ws.Range("B97").Value = 106.134217158065
$ws.Range("C97").Value = 120.746832746776

# Columns D-G: text values that look like numbers - force them to be stored
# as shared-string text (matching the rest of the sheet) rather than being
# auto-converted to numbers, then drop the temporary "Text" number format so
# the cell keeps the workbook's default (General) style, same as its peers.
foreach ($addr in @("D97", "E97", "F97", "G97")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D97").Value = "110.4"
$ws.Range("E97").Value = "112.4"
$ws.Range("F97").Value = " 88.4"
$ws.Range("G97").Value = "170.9"

foreach ($addr in @("D97", "E97", "F97", "G97")) {
    $ws.Range($addr).Style = "Normal"
}
